$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving its original style/format.
# Plain numeric-looking strings (e.g. "252.30") would otherwise be
# auto-converted to numbers by Excel, which would not match the source
# data (all Price/Volume columns in this sheet are text).
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" '42.914.43'
Set-TextValue "E2" '  +0.79%  '
Set-TextValue "D3" '2.288.10'
Set-TextValue "E3" '  +1.91%  '
Set-TextValue "E4" '  +0.16%  '
Set-TextValue "D5" '252.30'
Set-TextValue "E5" '  -0.01%  '
Set-TextValue "E6" '  +1.37%  '
Set-TextValue "D7" '73.96'
Set-TextValue "E7" '  +6.05%  '
Set-TextValue "E8" '  +0.14%  '
Set-TextValue "D9" '0.646'
Set-TextValue "E9" '  +1.15%  '
Set-TextValue "D10" '39.18'
Set-TextValue "E10" '  -4.51%  '
Set-TextValue "E11" '  +2.65%  '
Set-TextValue "D12" '59.05'
Set-TextValue "E12" '  -0.79%  '
Set-TextValue "E13" '  +1.86%  '
Set-TextValue "E14" '  +1.03%  '
Set-TextValue "D15" '2.631.09'
Set-TextValue "E15" '  +2.18%  '
Set-TextValue "D16" '15.32'
Set-TextValue "E16" '  +3.96%  '
Set-TextValue "E17" '  -2.06%  '
Set-TextValue "D18" '2.295.58'
Set-TextValue "E18" '  +2.17%  '
Set-TextValue "D19" '42.801.53'
Set-TextValue "E19" '  +1.03%  '
Set-TextValue "E20" '  +3.38%  '
Set-TextValue "E21" '  +1.31%  '
Set-TextValue "D22" '72.66'
Set-TextValue "E22" '  -0.27%  '
Set-TextValue "D23" '237.28'
Set-TextValue "E23" '  +1.11%  '
Set-TextValue "D24" '2.23'
Set-TextValue "E24" '  +6.96%  '
Set-TextValue "D25" '3.91'
Set-TextValue "E25" '  -1.62%  '
Set-TextValue "D26" '11.60'
Set-TextValue "E26" '  +0.21%  '
Set-TextValue "E27" '  -0.14%  '
Set-TextValue "E28" '  -0.34%  '
Set-TextValue "E29" '  -0.66%  '
Set-TextValue "D30" '2.19'
Set-TextValue "E30" '  -0.72%  '
Set-TextValue "D31" '167.15'
Set-TextValue "E31" '  -0.19%  '
Set-TextValue "D32" '21.05'
Set-TextValue "E32" '  +0.88%  '
Set-TextValue "E33" '  +6.19%  '
Set-TextValue "E34" '  +3.71%  '
Set-TextValue "D35" '0.0826'
Set-TextValue "E35" '  +4.96%  '
Set-TextValue "D36" '31.08'
Set-TextValue "E36" '  +11.52%  '
Set-TextValue "E37" '  +2.22%  '
Set-TextValue "E38" '  +12.05%  '
Set-TextValue "E39" '  +1.76%  '
Set-TextValue "E40" '  -2.73%  '
Set-TextValue "D41" '14.37'
Set-TextValue "E41" '  +14.26%  '
Set-TextValue "D42" '2.34'
Set-TextValue "E42" '  +3.26%  '
Set-TextValue "E43" '  +3.31%  '
Set-TextValue "E44" '  +7.91%  '
Set-TextValue "D45" '9.17'
Set-TextValue "E45" '  +4.83%  '
Set-TextValue "B46" 'MultiversX'
Set-TextValue "C46" 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue "D46" '61.87'
Set-TextValue "E46" '  -2.75%  '
Set-TextValue "B47" 'FTXToken'
Set-TextValue "C47" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D47" '4.88'
Set-TextValue "E47" '  -3.03%  '
Set-TextValue "E48" '  +1.81%  '
Set-TextValue "B49" 'BinanceUSD'
Set-TextValue "C49" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D49" '1.00'
Set-TextValue "E49" '  +0.26%  '
Set-TextValue "B50" 'ARBITRUM'
Set-TextValue "C50" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D50" '1.17'
Set-TextValue "E50" '  +0.07%  '
Set-TextValue "B51" 'Aave'
Set-TextValue "C51" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D51" '100.37'
Set-TextValue "E51" '  +6.15%  '
